$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I8").Value = "sv"
$ws.Range("J8").Value = "Statement-opinion"
$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"
$ws.Range("I13").Value = "sv"
$ws.Range("J13").Value = "Statement-opinion"
$ws.Range("I19").Value = "sd"
$ws.Range("J19").Value = "Statement-non-opinion"
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I23").Value = "sv"
$ws.Range("J23").Value = "Statement-opinion"
$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"
$ws.Range("I42").Value = "sv"
$ws.Range("J42").Value = "Statement-opinion"
$ws.Range("I45").Value = "sv"
$ws.Range("J45").Value = "Statement-opinion"
$ws.Range("I50").Value = "sd"
$ws.Range("J50").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "sv"
$ws.Range("J54").Value = "Statement-opinion"
$ws.Range("I55").Value = "sd"
$ws.Range("J55").Value = "Statement-non-opinion"
$ws.Range("I65").Value = "sd"
$ws.Range("J65").Value = "Statement-non-opinion"
$ws.Range("I81").Value = "b"
$ws.Range("J81").Value = "Acknowledge (Backchannel)"
$ws.Range("I93").Value = "sv"
$ws.Range("J93").Value = "Statement-opinion"
$ws.Range("I96").Value = "sv"
$ws.Range("J96").Value = "Statement-opinion"
$ws.Range("I106").Value = "sd"
$ws.Range("J106").Value = "Statement-non-opinion"
$ws.Range("I116").Value = "ba"
$ws.Range("J116").Value = "Appreciation"
$ws.Range("I123").Value = "aa"
$ws.Range("J123").Value = "Agree/Accept"
$ws.Range("I129").Value = "sd"
$ws.Range("J129").Value = "Statement-non-opinion"
$ws.Range("I149").Value = "sv"
$ws.Range("J149").Value = "Statement-opinion"
$ws.Range("I157").Value = "sd"
$ws.Range("J157").Value = "Statement-non-opinion"
$ws.Range("I159").Value = "ba"
$ws.Range("J159").Value = "Appreciation"
$ws.Range("I177").Value = "sd"
$ws.Range("J177").Value = "Statement-non-opinion"
$ws.Range("I178").Value = "%"
$ws.Range("J178").Value = "Uninterpretable"
$ws.Range("I183").Value = "sv"
$ws.Range("J183").Value = "Statement-opinion"
$ws.Range("I188").Value = "b"
$ws.Range("J188").Value = "Acknowledge (Backchannel)"
$ws.Range("I189").Value = "b"
$ws.Range("J189").Value = "Acknowledge (Backchannel)"
$ws.Range("I202").Value = "b"
$ws.Range("J202").Value = "Acknowledge (Backchannel)"
$ws.Range("I203").Value = "aa"
$ws.Range("J203").Value = "Agree/Accept"
$ws.Range("I207").Value = "sv"
$ws.Range("J207").Value = "Statement-opinion"
$ws.Range("I223").Value = "sd"
$ws.Range("J223").Value = "Statement-non-opinion"
$ws.Range("I225").Value = "sd"
$ws.Range("J225").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "sd"
$ws.Range("J243").Value = "Statement-non-opinion"
$ws.Range("I255").Value = "sd"
$ws.Range("J255").Value = "Statement-non-opinion"
$ws.Range("I274").Value = "sv"
$ws.Range("J274").Value = "Statement-opinion"
$ws.Range("I283").Value = "aa"
$ws.Range("J283").Value = "Agree/Accept"
$ws.Range("I285").Value = "sv"
$ws.Range("J285").Value = "Statement-opinion"
$ws.Range("I286").Value = "aa"
$ws.Range("J286").Value = "Agree/Accept"
$ws.Range("I287").Value = "sv"
$ws.Range("J287").Value = "Statement-opinion"
$ws.Range("I292").Value = "ba"
$ws.Range("J292").Value = "Appreciation"
